$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Add the two new header cells (G1:H1)
# ------------------------------------------------------------------
$ws.Range("G1").Value = "Table Name"
$ws.Range("H1").Value = "Value in table"

# ------------------------------------------------------------------
# 2. Add the "Table Name" column (G) for the rows that keep their
#    original position (rows 2-10 are unchanged in columns A-F)
# ------------------------------------------------------------------
$ws.Range("G2").Value  = "N/A"
$ws.Range("G3").Value  = "N/A"
$ws.Range("G4").Value  = "CEVAC_x_POWER_LATEST "
$ws.Range("H4").Value  = "Building Lighting actualvalue"
$ws.Range("G5").Value  = "CEVAC_x_TEMP_LATEST"
$ws.Range("H5").Value  = "All rooms"
$ws.Range("G6").Value  = "N/A"
$ws.Range("G7").Value  = "N/A"
$ws.Range("G8").Value  = "CEVAC_x_POWER_LATEST "
$ws.Range("H8").Value  = "sum of all actualvalue"
$ws.Range("G9").Value  = "CEVAC_x_POWER_LATEST "
$ws.Range("H9").Value  = "sum of all actualvalue"
$ws.Range("G10").Value = "CEVAC_x_POWER_LATEST "
$ws.Range("H10").Value = "sum of all actualvalue"

# ------------------------------------------------------------------
# 3. Rows 11-15 are reshuffled:
#      - old row 11 (Alert 10, "Any")        -> moves to row 13
#      - old row 12 (Alert 11, "IAQ/TBD")     -> removed entirely
#      - old row 13 (Alert 12, "Emergency")   -> moves to row 11
#      - old row 14 (Alert 13, "CO2" blank)   -> removed entirely
#      - old row 15 (Alert 14, "CO2")         -> moves to row 12
#    Simplest reliable way: wipe the old range and rewrite the three
#    surviving rows (with their new columns G/H) in the new order.
# ------------------------------------------------------------------
$ws.Range("A11:H15").ClearContents()

# New row 11 : Emergency (was row 13)
$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Emergency"
$ws.Range("C11").Value = "any"
$ws.Range("D11").Value = "any"
$ws.Range("E11").Value = ">80"
$ws.Range("G11").Value = "CEVAC_x_POWER_LATEST "
$ws.Range("H11").Value = "sum of 3rd floor emergency and basement emergency"

# New row 12 : CO2 (was row 15)
$ws.Range("A12").Value = 14
$ws.Range("B12").Value = "CO2"
$ws.Range("E12").Value = "1000 warn, 2000 alarm"
$ws.Range("G12").Value = "CEVAC_x_IAQ_LATEST"
$ws.Range("H12").Value = "all sensors"

# New row 13 : Any (was row 11)
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Any"
$ws.Range("C13").Value = "any"
$ws.Range("D13").Value = "any"
$ws.Range("E13").Value = "null/empty"
$ws.Range("G13").Value = "all"

# ------------------------------------------------------------------
# 4. Column G formatting - widen to fit its new contents
# ------------------------------------------------------------------
$ws.Columns("G").AutoFit()

# ------------------------------------------------------------------
# 5. Restore the active selection as recorded after the edit
# ------------------------------------------------------------------
$ws.Range("F17").Select()
